$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

$ws.Activate() | Out-Null
$ws.Range("K18").Select() | Out-Null
